# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-11-03 (45233) to 2023-11-13 (45243).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 101; $row++) {
    $ws.Cells.Item($row, 3).Value = 45243
}
